$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the test data that was populated in rows 1-28 (columns A:E).
# The rest of the used range (rows 29-542) was already empty.
# Re-applying the default style after clearing keeps the empty cell
# references present in the sheet (matching the original sparse layout)
# instead of Excel dropping now-blank rows/cells entirely.
$rng = $ws.Range("A1:E28")
$rng.Value = ""
$rng.Style = "Normal"
